$wb = $excel.ActiveWorkbook

$zhws = $wb.Worksheets.Item("zh-cn")
$dews = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4 and 5 originally shared the same handoff/handback
# timestamps (same shared-string entries), so both rows move forward together.
$zhws.Range("E4").Value = "2016-03-14 03:17:26"
$zhws.Range("H4").Value = "2016-03-14 03:17:40"
$zhws.Range("E5").Value = "2016-03-14 03:17:26"
$zhws.Range("H5").Value = "2016-03-14 03:17:40"

# de-de sheet: rows 4 and 5 originally shared the same handoff/handback
# timestamps (same shared-string entries), so both rows move forward together.
$dews.Range("E4").Value = "2016-03-14 03:17:30"
$dews.Range("H4").Value = "2016-03-14 03:17:46"
$dews.Range("E5").Value = "2016-03-14 03:17:30"
$dews.Range("H5").Value = "2016-03-14 03:17:46"
